# Financials update: insert a new "FY2018" data column at column D (period
# ending 2018-12-31, serial 43465), shifting the previously existing D:K
# columns to E:L. Five rows received a restated figure in the new column E
# (rather than simply inheriting the prior column D value verbatim), so
# those are corrected explicitly after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new blank column before column D; this pushes D:K -> E:L.
$ws.Columns("D").Insert()

# 2. The freshly inserted column has no formatting of its own; copy the
#    number formats (date / thousands styles) from the column immediately
#    to its right (the old column D, now shifted to E) so the new column
#    matches the rest of the table.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the new column D with the FY2018 figures.
$dValues = @{
    7=43465; 8=4727800; 9=2571300; 10=2156500; 12=291500; 13=0; 14=318500;
    15=0; 17=3931900; 18=795900; 20=82300; 21=1046400; 22=134500; 23=743700;
    24=80300; 25=0; 26=663400; 27=651100; 28=0; 29=-151900; 30=0; 31=0;
    32=-82300; 33=499200; 34=0; 35=499200; 38=43465; 41=161700; 42=0;
    43=2285200; 44=1097300; 45=486000; 46=4030200; 47=700; 48=1032600;
    49=4172400; 50=0; 51=0; 52=738400; 53=0; 54=9974300; 57=867500;
    58=547700; 59=1578200; 60=2993400; 61=2179000; 62=1591500; 63=0; 64=0;
    65=0; 66=6853200; 68=0; 69=0; 70=0; 71=0; 72=4334300; 73=0; 74=0; 75=0;
    76=3121100; 77=0; 80=43465; 81=499200; 83=168200; 84=0; 85=0; 86=0;
    87=0; 88=0; 89=368400; 91=-156600; 92=0; 93=0; 94=-130900; 96=-89200;
    97=0; 98=0; 99=0; 100=-363300; 101=4500; 102=-121300
}
foreach ($r in $dValues.Keys) {
    $ws.Range("D$r").Value = $dValues[$r]
}

# 4. A handful of rows were restated (their "after" column E value is not
#    simply the figure that used to live in column D) -- fix those up too.
$eValues = @{ 14=188600; 17=2594200; 18=284400; 20=-23600; 32=23600 }
foreach ($r in $eValues.Keys) {
    $ws.Range("E$r").Value = $eValues[$r]
}
